$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.358.28'
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('D3').Value = '1.845.57'
$ws.Range('E3').Value = '  -0.26%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9990'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6365'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07546'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.26%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2962'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.66'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.93%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07729'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.838.18'
$ws.Range('E12').Value = '  -1.35%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.991'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.74%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6832'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.57%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '83.09'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.21%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000009912'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.45%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.167'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.35%  '
$ws.Range('D18').Value = '29.394.64'
$ws.Range('E18').Value = '  -0.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '229.08'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.44'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.71%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.558'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.86%  '
$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('B24').Value = 'Monero'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '156.97'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.45%  '
$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1405'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.92%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.371'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.01%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.65'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.69%  '
$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.465'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.60%  '
$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.05704'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.16%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.249'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.13%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.129'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.032'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.73%  '
$ws.Range('B33').Value = 'LidoDAOToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.846'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.03%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.156'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.36%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7165'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.28%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.593'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.19%  '
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').Value = '1.257.60'
$ws.Range('E37').Value = '  +1.35%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.791'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.32%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01805'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.59%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9092'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.61%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.152'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.50%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.001'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.10%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.71'
$ws.Range('D43').Style = 'Normal'
$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').Value = '1.983.48'
$ws.Range('E44').Value = '  -2.11%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '66.49'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.41%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.196'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.45%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.047'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4025'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.704'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.21%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00000000116'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.79%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1127'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.57%  '
